$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C4").Value = "Domingo"
$ws.Range("D4").Value = "Pedro"
